$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel, since the source data must stay text.
foreach ($addr in @("D4","D5","D6","D7","D8","D17","D18","D19","D23","D25","D27","D29","D32","D34","D37","D38","D40","D41","D44","D48","D49","D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values from the latest cryptos data pull.
$ws.Range("D2").Value = '28.647.67'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '1.563.77'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '210.21'
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").Value = '0.512'
$ws.Range("E6").Value = '  +4.41%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '24.94'
$ws.Range("E8").Value = '  +5.61%  '
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("D12").Value = '1.788.45'
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").Value = '1.560.72'
$ws.Range("E13").Value = '  -0.16%  '
$ws.Range("D14").Value = '28.698.93'
$ws.Range("E14").Value = '  +1.44%  '
$ws.Range("E15").Value = '  +1.06%  '
$ws.Range("D17").Value = '61.48'
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").Value = '227.71'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '7.31'
$ws.Range("E19").Value = '  -0.43%  '
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = '9.02'
$ws.Range("E23").Value = '  +1.73%  '
$ws.Range("E24").Value = '  +1.61%  '
$ws.Range("D25").Value = '152.04'
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("E26").Value = '  +2.88%  '
$ws.Range("D27").Value = '14.77'
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("D29").Value = '6.22'
$ws.Range("E29").Value = '  -1.58%  '
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("D32").Value = '3.17'
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("D33").Value = '1.400.05'
$ws.Range("E33").Value = '  +1.72%  '
$ws.Range("D34").Value = '2.99'
$ws.Range("E34").Value = '  -3.04%  '
$ws.Range("E35").Value = '  -3.44%  '
$ws.Range("E36").Value = '  -1.93%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.30'
$ws.Range("E37").Value = '  -1.87%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.66'
$ws.Range("E38").Value = '  +0.98%  '
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '1.95'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("B41").Value = 'ImmutableX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D41").Value = '0.516'
$ws.Range("E41").Value = '  -0.63%  '
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").Value = '0.0459'
$ws.Range("E44").Value = '  -2.66%  '
$ws.Range("E46").Value = '  -1.79%  '
$ws.Range("D47").Value = '1.701.18'
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("D48").Value = '0.866'
$ws.Range("E48").Value = '  -5.65%  '
$ws.Range("D49").Value = '84.78'
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("D50").Value = '42.44'
$ws.Range("E50").Value = '  +5.94%  '
$ws.Range("E51").Value = '  -0.39%  '
